$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 97, shifting existing rows 97:157 down to 98:158
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new weekly record
$ws.Cells.Item(97, 1).Value = 4
$ws.Cells.Item(97, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(97, 3).Value = "Los Lagos"
$ws.Cells.Item(97, 4).Value = 44529
$ws.Cells.Item(97, 5).Value = 10
$ws.Cells.Item(97, 6).Value = 100112032
$ws.Cells.Item(97, 7).Value = "Zapallo italiano"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 80
$ws.Cells.Item(97, 11).Value = 10000
$ws.Cells.Item(97, 12).Value = 10000
$ws.Cells.Item(97, 13).Value = 10000
$ws.Cells.Item(97, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(97, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(97, 16).Value = 200
$ws.Cells.Item(97, 17).Value = 50
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# Copy the date cell style (s="2") from row 98 to the new row 97, column D
$ws.Range("D98").Copy()
$ws.Range("D97").PasteSpecial(-4122) | Out-Null
